# Generate Report for Handoff
# The 4f995e86-344a-492f-8d1f-3600612073de.md file has moved from
# "Handed back: in sync with en-US" to "Ready for handoff" for both
# locales, with refreshed handoff timestamps and a new handback-version
# warning recorded in the Error Detail column.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2d7a6590a347f14b80bd767ad6ef6efe39683b8/e2e/4f995e86-344a-492f-8d1f-3600612073de.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6ba725e815c8253b42abd832199010f48535d95/e2e/4f995e86-344a-492f-8d1f-3600612073de.md."

# --- Overview sheet: row 3 is the 4f995e86-...-.md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-07 09:14:30"

# --- zh-cn sheet: row 3 is the 4f995e86-...-.md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-07 09:14:25"
$wsZhCn.Range("P3").Value = $errorDetail
# Error Detail column widens to fit the new long message (stored width 40)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666666

# --- de-de sheet: row 3 is the 4f995e86-...-.md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-07 09:14:30"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666666
